$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.284.00'
$ws.Range("E2").Value = '  +0.49%  '

$ws.Range("D3").Value = '1.859.36'
$ws.Range("E3").Value = '  +0.38%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7039'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '238.46'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.33%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07991'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.17%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3029'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.53'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08197'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.67%  '

$ws.Range("D12").Value = '1.882.45'
$ws.Range("E12").Value = '  +1.60%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.203'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.04%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7080'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.33%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.74'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.79%  '

$ws.Range("D16").Value = '29.397.31'
$ws.Range("E16").Value = '  +0.86%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.842'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.80%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007914'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.31'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '238.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.18%  '

$ws.Range("D21").Value = '2.153.07'
$ws.Range("E21").Value = '  +2.35%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.12%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.484'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.21%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '163.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.81%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.897'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.73%  '

$ws.Range("E27").Value = '  +1.21%  '

$ws.Range("E28").Value = '  +0.28%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.926'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.08%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.422'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.76%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.478'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.48%  '

$ws.Range("E32").Value = '  -3.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.030'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.69%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05197'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.07%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.164'
$ws.Range("D35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7174'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.08%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.005'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.17%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.684'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.19%  '

$ws.Range("E39").Value = '  +0.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.729'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.10%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9453'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.53%  '

$ws.Range("D42").Value = '1.156.21'
$ws.Range("E42").Value = '  +5.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.007'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.05%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4272'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.78'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.38%  '

$ws.Range("E46").Value = '  +0.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.02'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5303'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.07%  '

$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.766'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.46%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.029.18'
$ws.Range("E50").Value = '  +1.57%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.172'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.23%  '
